# Obsidian vault auto-backup style edit: rebuild the "idéal percentage" block
# on sheet "20M ma gueule" as a vertical two-column table (F/G) instead of the
# old wide horizontal tables, and split the single "ideal %" goal-seek cell
# into two side-by-side goals (keep the 20M vs exceed the 10M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Row 3: the old rate header B3:U3 only keeps B3 - the rest are cleared
#    back to an empty (but still percent-styled) cell.
# ---------------------------------------------------------------------------
$ws.Range("C3:U3").ClearContents()

# ---------------------------------------------------------------------------
# 2) Drop the old "Formule" label that used to sit at F24, and the whole
#    horizontal PV table (F25:Y25) + its rate-label row (F26:Y26). Clear()
#    wipes both value AND formatting so the unused cells disappear entirely.
# ---------------------------------------------------------------------------
$ws.Range("F24").Clear()
$ws.Range("F25:Y25").Clear()
$ws.Range("F26:Y26").Clear()

# Re-apply the percent format to the two cells (F26:G26) that remain as
# empty, styled placeholders under the new table.
$ws.Range("F26:G26").NumberFormat = "0%"

# ---------------------------------------------------------------------------
# 3) Build the new vertical table in F6:G25 - column F holds the candidate
#    interest rate, column G the corresponding "-PV" (present value) for a
#    1 000 000 rente over 20 periods.
# ---------------------------------------------------------------------------
for ($r = 6; $r -le 25; $r++) {
    $rate = ($r - 5) / 100
    $cellF = $ws.Cells.Item($r, 6)
    $cellF.NumberFormat = "0%"
    $cellF.Value = $rate

    $ws.Cells.Item($r, 7).Formula = "=-1*PV(`$F$r,`$A`$25,`$B`$2)"
}

# ---------------------------------------------------------------------------
# 4) Update the two goal labels. F28 keeps referring to the same text slot
#    but the wording changes; G28 is a brand-new label for the second goal.
# ---------------------------------------------------------------------------
$ws.Range("F28").Value = "Pourcentage idéal pour garder les 20 M a la fin"
$ws.Range("G28").Value = "Pourcentage idéal pour dépasser les 10 M"

# ---------------------------------------------------------------------------
# 5) Row 29/30 used to hold a single goal-seek pair (F29/F30). Now there are
#    two: F-column keeps the "reach exactly 20 M" rate (near 0%), G-column
#    is the new "reach exactly 10 M" rate (what F29/F30 used to hold).
# ---------------------------------------------------------------------------
$ws.Range("G29").Formula = "=-1*PV(G30,`$A`$25,`$B`$2)"

$ws.Range("F30").Value = 0.0000000011815372845129261
$ws.Range("G30").NumberFormat = "0%"
$ws.Range("G30").Value = 0.077546895300093871

# ---------------------------------------------------------------------------
# 6) Column widths: F grows to fit the long label, G is new and sized to fit
#    its own (shorter) label.
# ---------------------------------------------------------------------------
$ws.Columns(6).ColumnWidth = 34.75
$ws.Columns(7).ColumnWidth = 31.666666666666668

# ---------------------------------------------------------------------------
# 7) Selection cosmetics, matching where the author's cursor ended up.
# ---------------------------------------------------------------------------
$ws.Range("F5").Select()
